# Femacal de La Calera - Haba: insert a new weekly record as row 159,
# shifting the existing rows 159-239 down to 160-240.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 159 (pushes old row 159..239 down to 160..240,
# carrying formatting/styles along, e.g. the date-format style on column D).
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159. The row mirrors the structure of its
# neighbours: Mercado/Region/Codreg/Categoria/Variedad/Calidad/Unidad/Region
# text columns stay constant for this subset, only the date & price/volume
# figures differ for this new record.
$ws.Range("A159").Value  = 3
$ws.Range("B159").Value  = "Femacal de La Calera"
$ws.Range("C159").Value  = "Coquimbo"
$ws.Range("D159").Value  = 44917
$ws.Range("E159").Value  = 5
$ws.Range("F159").Value  = 100112026
$ws.Range("G159").Value  = "Haba"
$ws.Range("H159").Value  = "Sin especificar"
$ws.Range("I159").Value  = "Primera"
$ws.Range("J159").Value  = 93
$ws.Range("K159").Value  = 8000
$ws.Range("L159").Value  = 8500
$ws.Range("M159").Value  = 8258
$ws.Range("N159").Value  = "$/saco 25 kilos"
$ws.Range("O159").Value  = "Provincia de Quillota"
$ws.Range("P159").Value  = 330
$ws.Range("Q159").Value  = 25
$ws.Range("R159").Value  = "Hortaliza"

# Make sure the date cell keeps the same number format as the rest of column D.
$ws.Range("D159").NumberFormat = $ws.Range("D160").NumberFormat
